$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("13866811204960522733", "19", "12", "2022", "Test", "3"),
    @("16767335235704918509", "11", "09", "2064", "Dennis 63th Bday", "5"),
    @("17531903181772362221", "13", "10", "2022", "TestEvent", "1"),
    @("13178649838081741293", "13", "10", "2022", "TestEvent2", "1")
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$c]
        $cell.Style = "Normal"
    }
}
